# Added Test for role based login
#
# 1. Insert a new "Login Details" sheet at the front of the workbook and a
#    new "Login Roles" sheet at the end, populate them with data, and
#    update the view/selection state on every sheet so that "Login Roles"
#    ends up the active tab.

$wb = $excel.ActiveWorkbook

# --- Create the two new sheets in the right slots --------------------------
# Add "Login Roles" first (placed before sheet 1 for the moment - we move it
# to the end once both new sheets exist), then add "Login Details" before
# sheet 1 too, so the order becomes:
#   Login Details, Create Opportunity, Create New Dealer, Login Roles
$loginRoles = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$loginRoles.Name = "Login Roles"

$loginDetails = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$loginDetails.Name = "Login Details"

# Move "Login Roles" (re-fetched by name, since the handle captured above can
# go stale once further sheets are inserted) to the very end.
$loginRolesWs = $wb.Worksheets.Item("Login Roles")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$loginRolesWs.Move([System.Reflection.Missing]::Value, $lastSheet)

# --- Populate "Login Details" ----------------------------------------------
$loginDetailsWs = $wb.Worksheets.Item("Login Details")

$loginDetailsWs.Range("A1").Value = "Username"
$loginDetailsWs.Range("B1").Value = "Password"
$loginDetailsWs.Range("A2").Value = "salestest1@cae.cae.qa"
$loginDetailsWs.Range("B2").Value = "CoxAuto123"

$loginDetailsWs.Range("A1:B2").WrapText = $true
$loginDetailsWs.Rows.Item(2).RowHeight = 28

# --- Populate "Login Roles" -------------------------------------------------
$loginRolesWs = $wb.Worksheets.Item("Login Roles")

$loginRolesWs.Columns.Item(1).ColumnWidth = 19.0
$loginRolesWs.Columns.Item(2).ColumnWidth = 13.83

$roleRows = @(
    @("User", "Role"),
    @("Watkins, Zanea", "Coordination Desk"),
    @("Reed, Carolyn", "Rewards Admin"),
    @("Graham, Darby", "All Sales, CAI Sales Ops"),
    @("Isom, Emilio", "All Sales, CAI Sales"),
    @("Levis, Ryan", "Reward Sales"),
    @("Abulafi, Iyad", "Sales Engineer, CAI Sales"),
    @("Jolitz, Maggie", "Business Admin:"),
    @("Lichtenberger, Lydia", "CAI Admin"),
    @("Assignee, Task", "CAI Chatter Only"),
    @("Augustaitis, George", "CAI Sales Ops Chatter"),
    @("Abrams, Randy", "CAI Sales Chatter"),
    @("Brunson, Wendy", "CoE Viewer (Chatter Plus)"),
    @("Augustine, Justin", "QA")
)

$tallRows = @(2, 4, 5, 7, 11, 13)

# Header row (User / Role) is written first, then the rest of column A
# (the "User" names) top-to-bottom, then the rest of column B (the "Role"
# names) top-to-bottom - this matches the shared-string allocation order
# baked into the target file (109=User, 110=Role, 111-123=names,
# 124-136=roles).
$loginRolesWs.Cells.Item(1, 1).Value = $roleRows[0][0]
$loginRolesWs.Cells.Item(1, 2).Value = $roleRows[0][1]

for ($i = 1; $i -lt $roleRows.Count; $i++) {
    $r = $i + 1
    $loginRolesWs.Cells.Item($r, 1).Value = $roleRows[$i][0]
}
for ($i = 1; $i -lt $roleRows.Count; $i++) {
    $r = $i + 1
    $loginRolesWs.Cells.Item($r, 2).Value = $roleRows[$i][1]
    $loginRolesWs.Cells.Item($r, 2).WrapText = $true
    if ($tallRows -contains $r) {
        $loginRolesWs.Rows.Item($r).RowHeight = 28
    }
}

# --- Update view/selection state on the two pre-existing sheets ------------
$createOpportunityWs = $wb.Worksheets.Item("Create Opportunity")
$createOpportunityWs.Range("A1:B2").Select()

$createNewDealerWs = $wb.Worksheets.Item("Create New Dealer")
$createNewDealerWs.Range("H2").Select()

# --- Select/activate the new sheets (Login Roles last so it becomes the
#     active tab, matching activeTab="3" / tabSelected="1") ----------------
$loginDetailsWs = $wb.Worksheets.Item("Login Details")
$loginDetailsWs.Range("D17").Select()

$loginRolesWs = $wb.Worksheets.Item("Login Roles")
$loginRolesWs.Range("C11").Select()
